# Adds a new "2022-Q4" sheet (right after "总计") with its fund-holdings
# table, and updates the "总计" (summary) sheet with the new quarter's
# row plus the resulting shift of all following rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Rebuild the "总计" summary sheet with the 2022-Q4 row inserted at
#    the top of the data (row 2), shifting everything else down by one.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$total.Cells.Clear()

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"
$totalHeader = $total.Range("B1:D1")
$totalHeader.Font.Bold = $true
$totalHeader.HorizontalAlignment = -4108
$totalHeader.VerticalAlignment = -4160
$totalHeader.Borders.LineStyle = 1

$totalRows = @(
    @("2022-Q4", 14, 4.83),
    @("2022-Q3", 19, 5.13),
    @("2021-Q4", 3, 0.16),
    @("2021-Q3", 1, 0.04),
    @("2020-Q4", 2, 0.26)
)

for ($i = 0; $i -lt $totalRows.Length; $i++) {
    $r = $i + 2
    $total.Cells.Item($r, 1).Value = $i
    $total.Cells.Item($r, 2).Value = $totalRows[$i][0]
    $total.Cells.Item($r, 3).Value = $totalRows[$i][1]
    $total.Cells.Item($r, 4).Value = $totalRows[$i][2]
}

$totalIdx = $total.Range("A2:A6")
$totalIdx.Font.Bold = $true
$totalIdx.HorizontalAlignment = -4108
$totalIdx.VerticalAlignment = -4160
$totalIdx.Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q4" worksheet right after "总计" and fill it
#    with the quarter's fund-holdings table.
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $total)
$q4.Name = "2022-Q4"

$q4Headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $q4Headers.Length; $i++) {
    $q4.Cells.Item(1, 2 + $i).Value = $q4Headers[$i]
}
$q4Header = $q4.Range("B1:H1")
$q4Header.Font.Bold = $true
$q4Header.HorizontalAlignment = -4108
$q4Header.VerticalAlignment = -4160
$q4Header.Borders.LineStyle = 1

# code, name, scale, stockPosition, positionRatio, marketValue, positionRank
$q4Funds = @(
    @("720001", "财通价值动量混合", "36.34", "79.52", "6.06", "2.2022", 3),
    @("001480", "财通成长优选混合", "20.59", "91.01", "5.96", "1.2272", 5),
    @("014915", "财通匠心优选一年持有期混合A", "5.48", "91.42", "6.26", "0.3430", 4),
    @("013142", "华商乐享互联灵活配置混合C", "5.21", "88.52", "3.30", "0.1719", 5),
    @("501046", "财通多策略福鑫定期开放灵活配置混合", "2.69", "91.79", "6.24", "0.1679", 3),
    @("008983", "财通科技创新混合A", "2.87", "87.06", "5.22", "0.1498", 5),
    @("001959", "华商乐享互联灵活配置混合A", "4.50", "88.52", "3.30", "0.1485", 5),
    @("009062", "财通智慧成长混合A", "2.05", "86.49", "5.85", "0.1199", 3),
    @("009063", "财通智慧成长混合C", "1.74", "86.49", "5.85", "0.1018", 3),
    @("008984", "财通科技创新混合C", "1.79", "87.06", "5.22", "0.0934", 5),
    @("014916", "财通匠心优选一年持有期混合C", "0.59", "91.42", "6.26", "0.0369", 4),
    @("001339", "兴银鼎新灵活配置混合", "0.71", "87.16", "4.72", "0.0335", 10),
    @("010124", "兴银景气优选混合A", "0.45", "83.23", "4.44", "0.0200", 10),
    @("010125", "兴银景气优选混合C", "0.35", "83.23", "4.44", "0.0155", 10)
)

# Columns B..G (fund code through market value) are stored as text in the
# source data (e.g. fund codes, and numbers kept at fixed precision such
# as "0.0200"), so force a text number-format before assigning the values
# to stop Excel from auto-coercing them to numeric cells.
$q4BG = $q4.Range("B2:G15")
$q4BG.NumberFormat = "@"

for ($i = 0; $i -lt $q4Funds.Length; $i++) {
    $r = $i + 2
    $q4.Cells.Item($r, 1).Value = $i
    $q4.Cells.Item($r, 2).Value = $q4Funds[$i][0]
    $q4.Cells.Item($r, 3).Value = $q4Funds[$i][1]
    $q4.Cells.Item($r, 4).Value = $q4Funds[$i][2]
    $q4.Cells.Item($r, 5).Value = $q4Funds[$i][3]
    $q4.Cells.Item($r, 6).Value = $q4Funds[$i][4]
    $q4.Cells.Item($r, 7).Value = $q4Funds[$i][5]
    $q4.Cells.Item($r, 8).Value = $q4Funds[$i][6]
}

$q4Idx = $q4.Range("A2:A15")
$q4Idx.Font.Bold = $true
$q4Idx.HorizontalAlignment = -4108
$q4Idx.VerticalAlignment = -4160
$q4Idx.Borders.LineStyle = 1

# Keep the originally-selected tab ("2020-Q4", the last sheet) active,
# since adding the new sheet would otherwise steal the selection.
$wb.Worksheets.Item("2020-Q4").Activate()

Write-Output "done"
